$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (id, speaker_variant) values per row, 2..18
$data = @{
    2  = @("#griet", "Griet")
    3  = @("#mr", "Mr")
    4  = @("#blinde", "Blinde")
    5  = @("#modde", "Modde")
    6  = @("#blind:", "Blind:")
    7  = @("#reyn", "Reyn")
    8  = @("#schout", "Schout")
    9  = @("#ian", "Ian")
    10 = @("#reyn:", "Reyn:")
    11 = @("#kreup", "Kreup")
    12 = @("#wijff", "wijff")
    13 = @("#kreup:", "Kreup:")
    14 = @("#moer", "Moer")
    15 = @("#mich:", "Mich:")
    16 = @("#wijff", "Wijff")
    17 = @("#melis", "Melis")
    18 = @("#lijs", "Lijs")
}

foreach ($row in 2..18) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $null
}
